$wb = $excel.ActiveWorkbook

# --- Helper: apply header styling (bold, thin border all sides, center/top align) ---
function Set-HeaderStyle($rng) {
    $rng.Font.Bold = $true
    $rng.Borders.LineStyle = 1
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
}

# --- 1) Create the two brand-new sheets first. Worksheet handles in this ---
#        engine are position-bound, so we avoid reusing any worksheet ---
#        object across further Add()/Move() calls - we only rename here, ---
#        then re-fetch everything fresh by name afterwards. ---
$newSheet1 = $wb.Worksheets.Add()
$newSheet1.Name = "Player Info"
$newSheet2 = $wb.Worksheets.Add()
$newSheet2.Name = "ODI Batting Extra"

# --- 2) Re-order sheets: Player Info, ODI Batting, ODI Bowling, ODI Batting Extra ---
$playerInfo = $wb.Worksheets.Item("Player Info")
$playerInfo.Move($wb.Worksheets.Item(1))

$lastIdx = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIdx)
$extraSheet = $wb.Worksheets.Item("ODI Batting Extra")
$extraSheet.Move($null, $lastSheet)

# Sanity: expected order now is Player Info, ODI Batting, ODI Bowling, ODI Batting Extra

# --- 3) Fetch fresh (stable) worksheet handles by name for all data edits below ---
$playerInfo = $wb.Worksheets.Item("Player Info")
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBowling = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Item("ODI Batting Extra")

# --- 4) Populate "Player Info" sheet ---
$playerInfo.Cells.Item(1,1).Value = "ID"
$playerInfo.Cells.Item(1,2).Value = "NAME"
$playerInfo.Cells.Item(1,3).Value = "BATTING_HAND"
$playerInfo.Cells.Item(1,4).Value = "BOWL_STYLE"
Set-HeaderStyle($playerInfo.Range("A1:D1"))

$playerInfo.Cells.Item(2,1).NumberFormat = "@"
$playerInfo.Cells.Item(2,1).Value = "4268"
$playerInfo.Cells.Item(2,2).Value = "Sabbir Rahman"
$playerInfo.Cells.Item(2,3).Value = "Right Handed"
$playerInfo.Cells.Item(2,4).Value = "Right Arm Leg Break"

# --- 5) Update "ODI Batting" sheet: rename MATCH_CARD_LINK -> MATCH_CODE, ---
#        replace URL values in column D with bare match codes, ---
#        and clear the stray empty INNING_NUMBER cells ---
$odiBatting.Cells.Item(1,4).Value = "MATCH_CODE"
$odiBatting.Cells.Item(2,4).NumberFormat = "@"
$odiBatting.Cells.Item(2,4).Value = "3699"
$odiBatting.Cells.Item(3,4).NumberFormat = "@"
$odiBatting.Cells.Item(3,4).Value = "3701"
$odiBatting.Cells.Item(4,4).NumberFormat = "@"
$odiBatting.Cells.Item(4,4).Value = "3702"
$odiBatting.Cells.Item(5,4).NumberFormat = "@"
$odiBatting.Cells.Item(5,4).Value = "3704"
$odiBatting.Cells.Item(6,4).NumberFormat = "@"
$odiBatting.Cells.Item(6,4).Value = "3708"
$odiBatting.Cells.Item(7,4).NumberFormat = "@"
$odiBatting.Cells.Item(7,4).Value = "3754"
$odiBatting.Cells.Item(8,4).NumberFormat = "@"
$odiBatting.Cells.Item(8,4).Value = "3765"
$odiBatting.Cells.Item(9,4).NumberFormat = "@"
$odiBatting.Cells.Item(9,4).Value = "3774"
$odiBatting.Cells.Item(10,4).NumberFormat = "@"
$odiBatting.Cells.Item(10,4).Value = "3780"
$odiBatting.Cells.Item(11,4).NumberFormat = "@"
$odiBatting.Cells.Item(11,4).Value = "3784"
$odiBatting.Cells.Item(12,4).NumberFormat = "@"
$odiBatting.Cells.Item(12,4).Value = "3791"
$odiBatting.Cells.Item(13,4).NumberFormat = "@"
$odiBatting.Cells.Item(13,4).Value = "3797"
$odiBatting.Cells.Item(14,4).NumberFormat = "@"
$odiBatting.Cells.Item(14,4).Value = "3798"
$odiBatting.Cells.Item(15,4).NumberFormat = "@"
$odiBatting.Cells.Item(15,4).Value = "3799"
$odiBatting.Cells.Item(16,4).NumberFormat = "@"
$odiBatting.Cells.Item(16,4).Value = "3808"
$odiBatting.Cells.Item(17,4).NumberFormat = "@"
$odiBatting.Cells.Item(17,4).Value = "3810"
$odiBatting.Cells.Item(18,4).NumberFormat = "@"
$odiBatting.Cells.Item(18,4).Value = "3811"
$odiBatting.Cells.Item(19,4).NumberFormat = "@"
$odiBatting.Cells.Item(19,4).Value = "3813"
$odiBatting.Cells.Item(20,4).NumberFormat = "@"
$odiBatting.Cells.Item(20,4).Value = "3816"
$odiBatting.Cells.Item(21,4).NumberFormat = "@"
$odiBatting.Cells.Item(21,4).Value = "3818"
$odiBatting.Cells.Item(22,4).NumberFormat = "@"
$odiBatting.Cells.Item(22,4).Value = "3854"
$odiBatting.Cells.Item(23,4).NumberFormat = "@"
$odiBatting.Cells.Item(23,4).Value = "3856"
$odiBatting.Cells.Item(24,4).NumberFormat = "@"
$odiBatting.Cells.Item(24,4).Value = "3857"
$odiBatting.Cells.Item(25,4).NumberFormat = "@"
$odiBatting.Cells.Item(25,4).Value = "3936"
$odiBatting.Cells.Item(26,4).NumberFormat = "@"
$odiBatting.Cells.Item(26,4).Value = "3938"
$odiBatting.Cells.Item(27,4).NumberFormat = "@"
$odiBatting.Cells.Item(27,4).Value = "3941"
$odiBatting.Cells.Item(28,4).NumberFormat = "@"
$odiBatting.Cells.Item(28,4).Value = "3946"
$odiBatting.Cells.Item(29,4).NumberFormat = "@"
$odiBatting.Cells.Item(29,4).Value = "3948"
$odiBatting.Cells.Item(30,4).NumberFormat = "@"
$odiBatting.Cells.Item(30,4).Value = "3949"
$odiBatting.Cells.Item(31,4).NumberFormat = "@"
$odiBatting.Cells.Item(31,4).Value = "3969"
$odiBatting.Cells.Item(32,4).NumberFormat = "@"
$odiBatting.Cells.Item(32,4).Value = "3970"
$odiBatting.Cells.Item(33,4).NumberFormat = "@"
$odiBatting.Cells.Item(33,4).Value = "3971"
$odiBatting.Cells.Item(34,4).NumberFormat = "@"
$odiBatting.Cells.Item(34,4).Value = "4011"
$odiBatting.Cells.Item(35,4).NumberFormat = "@"
$odiBatting.Cells.Item(35,4).Value = "4012"
$odiBatting.Cells.Item(36,4).NumberFormat = "@"
$odiBatting.Cells.Item(36,4).Value = "4014"
$odiBatting.Cells.Item(37,4).NumberFormat = "@"
$odiBatting.Cells.Item(37,4).Value = "4022"
$odiBatting.Cells.Item(38,4).NumberFormat = "@"
$odiBatting.Cells.Item(38,4).Value = "4024"
$odiBatting.Cells.Item(39,4).NumberFormat = "@"
$odiBatting.Cells.Item(39,4).Value = "4025"
$odiBatting.Cells.Item(40,4).NumberFormat = "@"
$odiBatting.Cells.Item(40,4).Value = "4027"
$odiBatting.Cells.Item(41,4).NumberFormat = "@"
$odiBatting.Cells.Item(41,4).Value = "4031"
$odiBatting.Cells.Item(42,4).NumberFormat = "@"
$odiBatting.Cells.Item(42,4).Value = "4035"
$odiBatting.Cells.Item(43,4).NumberFormat = "@"
$odiBatting.Cells.Item(43,4).Value = "4039"
$odiBatting.Cells.Item(44,4).NumberFormat = "@"
$odiBatting.Cells.Item(44,4).Value = "4047"
$odiBatting.Cells.Item(45,4).NumberFormat = "@"
$odiBatting.Cells.Item(45,4).Value = "4080"
$odiBatting.Cells.Item(46,4).NumberFormat = "@"
$odiBatting.Cells.Item(46,4).Value = "4083"
$odiBatting.Cells.Item(47,4).NumberFormat = "@"
$odiBatting.Cells.Item(47,4).Value = "4086"
$odiBatting.Cells.Item(48,4).NumberFormat = "@"
$odiBatting.Cells.Item(48,4).Value = "4109"
$odiBatting.Cells.Item(49,4).NumberFormat = "@"
$odiBatting.Cells.Item(49,4).Value = "4116"
$odiBatting.Cells.Item(50,4).NumberFormat = "@"
$odiBatting.Cells.Item(50,4).Value = "4121"
$odiBatting.Cells.Item(51,4).NumberFormat = "@"
$odiBatting.Cells.Item(51,4).Value = "4122"
$odiBatting.Cells.Item(52,4).NumberFormat = "@"
$odiBatting.Cells.Item(52,4).Value = "4124"
$odiBatting.Cells.Item(53,4).NumberFormat = "@"
$odiBatting.Cells.Item(53,4).Value = "4179"
$odiBatting.Cells.Item(54,4).NumberFormat = "@"
$odiBatting.Cells.Item(54,4).Value = "4180"
$odiBatting.Cells.Item(55,4).NumberFormat = "@"
$odiBatting.Cells.Item(55,4).Value = "4181"
$odiBatting.Cells.Item(56,4).NumberFormat = "@"
$odiBatting.Cells.Item(56,4).Value = "4250"
$odiBatting.Cells.Item(57,4).NumberFormat = "@"
$odiBatting.Cells.Item(57,4).Value = "4251"
$odiBatting.Cells.Item(58,4).NumberFormat = "@"
$odiBatting.Cells.Item(58,4).Value = "4252"
$odiBatting.Cells.Item(59,4).NumberFormat = "@"
$odiBatting.Cells.Item(59,4).Value = "4286"
$odiBatting.Cells.Item(60,4).NumberFormat = "@"
$odiBatting.Cells.Item(60,4).Value = "4293"
$odiBatting.Cells.Item(61,4).NumberFormat = "@"
$odiBatting.Cells.Item(61,4).Value = "4295"
$odiBatting.Cells.Item(62,4).NumberFormat = "@"
$odiBatting.Cells.Item(62,4).Value = "4296"
$odiBatting.Cells.Item(63,4).NumberFormat = "@"
$odiBatting.Cells.Item(63,4).Value = "4329"
$odiBatting.Cells.Item(64,4).NumberFormat = "@"
$odiBatting.Cells.Item(64,4).Value = "4345"
$odiBatting.Cells.Item(65,4).NumberFormat = "@"
$odiBatting.Cells.Item(65,4).Value = "4356"
$odiBatting.Cells.Item(66,4).NumberFormat = "@"
$odiBatting.Cells.Item(66,4).Value = "4357"
$odiBatting.Cells.Item(67,4).NumberFormat = "@"
$odiBatting.Cells.Item(67,4).Value = "4358"

# Clear stray empty INNING_NUMBER cells (did not bat rows)
$odiBatting.Cells.Item(14,2).ClearContents()
$odiBatting.Cells.Item(15,2).ClearContents()
$odiBatting.Cells.Item(20,2).ClearContents()
$odiBatting.Cells.Item(21,2).ClearContents()
$odiBatting.Cells.Item(35,2).ClearContents()
$odiBatting.Cells.Item(48,2).ClearContents()
$odiBatting.Cells.Item(59,2).ClearContents()

# --- 6) Update "ODI Bowling" sheet: rename MATCH_CARD_LINK -> MATCH_CODE, ---
#        replace URL values in column B with bare match codes ---
$odiBowling.Cells.Item(1,2).Value = "MATCH_CODE"
$odiBowling.Cells.Item(2,2).NumberFormat = "@"
$odiBowling.Cells.Item(2,2).Value = "3699"
$odiBowling.Cells.Item(3,2).NumberFormat = "@"
$odiBowling.Cells.Item(3,2).Value = "3701"
$odiBowling.Cells.Item(4,2).NumberFormat = "@"
$odiBowling.Cells.Item(4,2).Value = "3702"
$odiBowling.Cells.Item(5,2).NumberFormat = "@"
$odiBowling.Cells.Item(5,2).Value = "3754"
$odiBowling.Cells.Item(6,2).NumberFormat = "@"
$odiBowling.Cells.Item(6,2).Value = "3765"
$odiBowling.Cells.Item(7,2).NumberFormat = "@"
$odiBowling.Cells.Item(7,2).Value = "3774"
$odiBowling.Cells.Item(8,2).NumberFormat = "@"
$odiBowling.Cells.Item(8,2).Value = "3780"
$odiBowling.Cells.Item(9,2).NumberFormat = "@"
$odiBowling.Cells.Item(9,2).Value = "3784"
$odiBowling.Cells.Item(10,2).NumberFormat = "@"
$odiBowling.Cells.Item(10,2).Value = "3799"
$odiBowling.Cells.Item(11,2).NumberFormat = "@"
$odiBowling.Cells.Item(11,2).Value = "3813"
$odiBowling.Cells.Item(12,2).NumberFormat = "@"
$odiBowling.Cells.Item(12,2).Value = "3857"
$odiBowling.Cells.Item(13,2).NumberFormat = "@"
$odiBowling.Cells.Item(13,2).Value = "3938"
$odiBowling.Cells.Item(14,2).NumberFormat = "@"
$odiBowling.Cells.Item(14,2).Value = "3971"
$odiBowling.Cells.Item(15,2).NumberFormat = "@"
$odiBowling.Cells.Item(15,2).Value = "4031"
$odiBowling.Cells.Item(16,2).NumberFormat = "@"
$odiBowling.Cells.Item(16,2).Value = "4047"
$odiBowling.Cells.Item(17,2).NumberFormat = "@"
$odiBowling.Cells.Item(17,2).Value = "4083"
$odiBowling.Cells.Item(18,2).NumberFormat = "@"
$odiBowling.Cells.Item(18,2).Value = "4086"
$odiBowling.Cells.Item(19,2).NumberFormat = "@"
$odiBowling.Cells.Item(19,2).Value = "4250"
$odiBowling.Cells.Item(20,2).NumberFormat = "@"
$odiBowling.Cells.Item(20,2).Value = "4251"
$odiBowling.Cells.Item(21,2).NumberFormat = "@"
$odiBowling.Cells.Item(21,2).Value = "4296"
$odiBowling.Cells.Item(22,2).NumberFormat = "@"
$odiBowling.Cells.Item(22,2).Value = "4357"

# --- 7) Populate "ODI Batting Extra" sheet ---
$extra.Cells.Item(1,1).Value = "MATCH_CODE"
$extra.Cells.Item(1,2).Value = "BATTING_POSITION"
$extra.Cells.Item(1,3).Value = "NUM_4"
$extra.Cells.Item(1,4).Value = "NUM_6"
$extra.Cells.Item(1,5).Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Cells.Item(1,6).Value = "MAN_OF_MATCH"
Set-HeaderStyle($extra.Range("A1:F1"))

$extra.Cells.Item(2,1).NumberFormat = "@"
$extra.Cells.Item(2,1).Value = "4109"
$extra.Cells.Item(2,2).Value = 7
$extra.Cells.Item(2,6).Value = "NO"

$extra.Cells.Item(3,1).NumberFormat = "@"
$extra.Cells.Item(3,1).Value = "4116"
$extra.Cells.Item(3,2).Value = 6
$extra.Cells.Item(3,3).NumberFormat = "@"
$extra.Cells.Item(3,3).Value = "3"
$extra.Cells.Item(3,4).NumberFormat = "@"
$extra.Cells.Item(3,4).Value = "1"
$extra.Cells.Item(3,5).NumberFormat = "@"
$extra.Cells.Item(3,5).Value = "7.50%"
$extra.Cells.Item(3,6).Value = "NO"

$extra.Cells.Item(4,1).NumberFormat = "@"
$extra.Cells.Item(4,1).Value = "4121"
$extra.Cells.Item(4,2).Value = 6
$extra.Cells.Item(4,3).NumberFormat = "@"
$extra.Cells.Item(4,3).Value = "0"
$extra.Cells.Item(4,4).NumberFormat = "@"
$extra.Cells.Item(4,4).Value = "0"
$extra.Cells.Item(4,5).NumberFormat = "@"
$extra.Cells.Item(4,5).Value = "2.78%"
$extra.Cells.Item(4,6).Value = "NO"

$extra.Cells.Item(5,1).NumberFormat = "@"
$extra.Cells.Item(5,1).Value = "4122"
$extra.Cells.Item(5,2).Value = 6
$extra.Cells.Item(5,3).NumberFormat = "@"
$extra.Cells.Item(5,3).Value = "2"
$extra.Cells.Item(5,4).NumberFormat = "@"
$extra.Cells.Item(5,4).Value = "0"
$extra.Cells.Item(5,5).NumberFormat = "@"
$extra.Cells.Item(5,5).Value = "12.20%"
$extra.Cells.Item(5,6).Value = "NO"

$extra.Cells.Item(6,1).NumberFormat = "@"
$extra.Cells.Item(6,1).Value = "4124"
$extra.Cells.Item(6,2).Value = 3
$extra.Cells.Item(6,3).NumberFormat = "@"
$extra.Cells.Item(6,3).Value = "0"
$extra.Cells.Item(6,4).NumberFormat = "@"
$extra.Cells.Item(6,4).Value = "0"
$extra.Cells.Item(6,5).NumberFormat = "@"
$extra.Cells.Item(6,5).Value = "1.41%"
$extra.Cells.Item(6,6).Value = "NO"

$extra.Cells.Item(7,1).NumberFormat = "@"
$extra.Cells.Item(7,1).Value = "4179"
$extra.Cells.Item(7,6).Value = "NO"

$extra.Cells.Item(8,1).NumberFormat = "@"
$extra.Cells.Item(8,1).Value = "4180"
$extra.Cells.Item(8,2).Value = 6
$extra.Cells.Item(8,3).NumberFormat = "@"
$extra.Cells.Item(8,3).Value = "1"
$extra.Cells.Item(8,4).NumberFormat = "@"
$extra.Cells.Item(8,4).Value = "0"
$extra.Cells.Item(8,5).NumberFormat = "@"
$extra.Cells.Item(8,5).Value = "4.48%"
$extra.Cells.Item(8,6).Value = "NO"

$extra.Cells.Item(9,1).NumberFormat = "@"
$extra.Cells.Item(9,1).Value = "4181"
$extra.Cells.Item(9,2).Value = 7
$extra.Cells.Item(9,3).NumberFormat = "@"
$extra.Cells.Item(9,3).Value = "2"
$extra.Cells.Item(9,4).NumberFormat = "@"
$extra.Cells.Item(9,4).Value = "0"
$extra.Cells.Item(9,5).NumberFormat = "@"
$extra.Cells.Item(9,5).Value = "3.99%"
$extra.Cells.Item(9,6).Value = "NO"

$extra.Cells.Item(10,1).NumberFormat = "@"
$extra.Cells.Item(10,1).Value = "4250"
$extra.Cells.Item(10,2).Value = 7
$extra.Cells.Item(10,3).NumberFormat = "@"
$extra.Cells.Item(10,3).Value = "2"
$extra.Cells.Item(10,4).NumberFormat = "@"
$extra.Cells.Item(10,4).Value = "0"
$extra.Cells.Item(10,5).NumberFormat = "@"
$extra.Cells.Item(10,5).Value = "5.60%"
$extra.Cells.Item(10,6).Value = "NO"

$extra.Cells.Item(11,1).NumberFormat = "@"
$extra.Cells.Item(11,1).Value = "4251"
$extra.Cells.Item(11,2).Value = 7
$extra.Cells.Item(11,3).NumberFormat = "@"
$extra.Cells.Item(11,3).Value = "7"
$extra.Cells.Item(11,4).NumberFormat = "@"
$extra.Cells.Item(11,4).Value = "0"
$extra.Cells.Item(11,5).NumberFormat = "@"
$extra.Cells.Item(11,5).Value = "19.03%"
$extra.Cells.Item(11,6).Value = "NO"

$extra.Cells.Item(12,1).NumberFormat = "@"
$extra.Cells.Item(12,1).Value = "4252"
$extra.Cells.Item(12,2).Value = 6
$extra.Cells.Item(12,3).NumberFormat = "@"
$extra.Cells.Item(12,3).Value = "12"
$extra.Cells.Item(12,4).NumberFormat = "@"
$extra.Cells.Item(12,4).Value = "2"
$extra.Cells.Item(12,5).NumberFormat = "@"
$extra.Cells.Item(12,5).Value = "42.15%"
$extra.Cells.Item(12,6).Value = "NO"

$extra.Cells.Item(13,1).NumberFormat = "@"
$extra.Cells.Item(13,1).Value = "4286"
$extra.Cells.Item(13,2).Value = 7
$extra.Cells.Item(13,6).Value = "NO"

$extra.Cells.Item(14,1).NumberFormat = "@"
$extra.Cells.Item(14,1).Value = "4293"
$extra.Cells.Item(14,2).Value = 7
$extra.Cells.Item(14,3).NumberFormat = "@"
$extra.Cells.Item(14,3).Value = "0"
$extra.Cells.Item(14,4).NumberFormat = "@"
$extra.Cells.Item(14,4).Value = "0"
$extra.Cells.Item(14,6).Value = "NO"

$extra.Cells.Item(15,1).NumberFormat = "@"
$extra.Cells.Item(15,1).Value = "4295"
$extra.Cells.Item(15,6).Value = "NO"

$extra.Cells.Item(16,1).NumberFormat = "@"
$extra.Cells.Item(16,1).Value = "4296"
$extra.Cells.Item(16,2).Value = 3
$extra.Cells.Item(16,3).NumberFormat = "@"
$extra.Cells.Item(16,3).Value = "0"
$extra.Cells.Item(16,4).NumberFormat = "@"
$extra.Cells.Item(16,4).Value = "0"
$extra.Cells.Item(16,6).Value = "NO"

$extra.Cells.Item(17,1).NumberFormat = "@"
$extra.Cells.Item(17,1).Value = "4329"
$extra.Cells.Item(17,6).Value = "NO"

$extra.Cells.Item(18,1).NumberFormat = "@"
$extra.Cells.Item(18,1).Value = "4345"
$extra.Cells.Item(18,2).Value = 7
$extra.Cells.Item(18,3).NumberFormat = "@"
$extra.Cells.Item(18,3).Value = "5"
$extra.Cells.Item(18,4).NumberFormat = "@"
$extra.Cells.Item(18,4).Value = "0"
$extra.Cells.Item(18,5).NumberFormat = "@"
$extra.Cells.Item(18,5).Value = "12.59%"
$extra.Cells.Item(18,6).Value = "NO"

$extra.Cells.Item(19,1).NumberFormat = "@"
$extra.Cells.Item(19,1).Value = "4356"
$extra.Cells.Item(19,6).Value = "NO"

$extra.Cells.Item(20,1).NumberFormat = "@"
$extra.Cells.Item(20,1).Value = "4357"
$extra.Cells.Item(20,6).Value = "NO"

$extra.Cells.Item(21,1).NumberFormat = "@"
$extra.Cells.Item(21,1).Value = "4358"
$extra.Cells.Item(21,6).Value = "NO"

